$d = $word.ActiveDocument
$find = $d.Content.Find

function Replace-Text($oldText, $newText) {
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# M2DocEvaluator.caseQuery : 559 -> 586
Replace-Text "M2DocEvaluator.java:559)" "M2DocEvaluator.java:586)"

# M2DocEvaluator.doSwitch : 1216 -> 1239 (3 occurrences, replaced in one pass)
Replace-Text "M2DocEvaluator.java:1216)" "M2DocEvaluator.java:1239)"

# M2DocEvaluator.caseBlock : 1425 -> 1464
Replace-Text "M2DocEvaluator.java:1425)" "M2DocEvaluator.java:1464)"

# M2DocEvaluator.caseDocumentTemplate : 287 -> 296
Replace-Text "M2DocEvaluator.java:287)" "M2DocEvaluator.java:296)"

# M2DocEvaluator.generate : 276 -> 281
Replace-Text "M2DocEvaluator.java:276)" "M2DocEvaluator.java:281)"

# M2DocUtils.generate : 694 -> 805
Replace-Text "M2DocUtils.java:694)" "M2DocUtils.java:805)"

# AbstractTemplatesTestSuite.prepareoutputAndGenerate : 480 -> 511
Replace-Text "AbstractTemplatesTestSuite.java:480)" "AbstractTemplatesTestSuite.java:511)"

# AbstractTemplatesTestSuite.generation : 389 -> 420
Replace-Text "AbstractTemplatesTestSuite.java:389)" "AbstractTemplatesTestSuite.java:420)"

# Insert a new stack frame line before the second occurrence of
# "RunAfters.evaluate(RunAfters.java:27)" (the one that directly follows
# "ParentRunner$2.evaluate(ParentRunner.java:268)"), which is the only
# unique anchor for that particular occurrence.
$tab = [char]9
$lf = [char]10
$anchor = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $lf + $tab + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
$replacement = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $lf + $tab + "at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + $lf + $tab + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
Replace-Text $anchor $replacement
